{"js": "// Replace the supervisor's signature name \"\u041c. \u0421. \u0410\u043d\u0430\u043d\u044c\u0435\u0432\u0441\u043a\u0438\u0439\" with\n// \"\u0415. \u041c. \u041a\u0443\u0437\u043d\u0435\u0446\u043e\u0432\u0430\" on the \"\u0420\u0443\u043a\u043e\u0432\u043e\u0434\u0438\u0442\u0435\u043b\u044c ... ____  \u041c. \u0421. \u0410\u043d\u0430\u043d\u044c\u0435\u0432\u0441\u043a\u0438\u0439\" line.\n//\n// The target run (\" \u041c. \u0421. \u0410\u043d\u0430\u043d\u044c\u0435\u0432\u0441\u043a\u0438\u0439\") is split, in the canonical edit,\n// into six sibling runs instead of being collapsed into one:\n//   \" \" | \"\u0415\" | \". \" | \"\u041c\" | \". \" | \"\u041a\u0443\u0437\u043d\u0435\u0446\u043e\u0432\u0430\"\n// A plain insertText(\"...\", \"Replace\") would leave the new text in a single\n// run (and the host even coalesces it into the preceding \"___\" run), so we\n// rebuild that exact run sequence with insertOoxml (a \"flat OPC\" package),\n// which inserts the literal OOXML runs verbatim without merging them.\n\nconst body = context.document.body;\nconst hits = body.search(\" \u041c. \u0421. \u0410\u043d\u0430\u043d\u044c\u0435\u0432\u0441\u043a\u0438\u0439\", { matchCase: true });\nhits.load(\"items\");\nawait context.sync();\n\nif (hits.items.length === 0) {\n  throw new Error('Target text \" \u041c. \u0421. \u0410\u043d\u0430\u043d\u044c\u0435\u0432\u0441\u043a\u0438\u0439\" was not found in the document body.');\n}\n\nconst runProps = '<w:rPr><w:rFonts w:cs=\"Times New Roman\"/><w:szCs w:val=\"28\"/></w:rPr>';\nconst flatOpc =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\"><pkg:xmlData>' +\n      '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n        '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n      '</Relationships>' +\n    '</pkg:xmlData></pkg:part>' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n      '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p>' +\n        '<w:r>' + runProps + '<w:t xml:space=\"preserve\"> </w:t></w:r>' +\n        '<w:r>' + runProps + '<w:t>\u0415</w:t></w:r>' +\n        '<w:r>' + runProps + '<w:t xml:space=\"preserve\">. </w:t></w:r>' +\n        '<w:r>' + runProps + '<w:t>\u041c</w:t></w:r>' +\n        '<w:r>' + runProps + '<w:t xml:space=\"preserve\">. </w:t></w:r>' +\n        '<w:r>' + runProps + '<w:t>\u041a\u0443\u0437\u043d\u0435\u0446\u043e\u0432\u0430</w:t></w:r>' +\n      '</w:p></w:body></w:document>' +\n    '</pkg:xmlData></pkg:part>' +\n  '</pkg:package>';\n\nhits.items[0].insertOoxml(flatOpc, \"Replace\");\nawait context.sync();\n", "ps1": "# Replace the supervisor's signature name \"\u041c. \u0421. \u0410\u043d\u0430\u043d\u044c\u0435\u0432\u0441\u043a\u0438\u0439\" with\n# \"\u0415. \u041c. \u041a\u0443\u0437\u043d\u0435\u0446\u043e\u0432\u0430\" on the \"\u0420\u0443\u043a\u043e\u0432\u043e\u0434\u0438\u0442\u0435\u043b\u044c ... ____  \u041c. \u0421. \u0410\u043d\u0430\u043d\u044c\u0435\u0432\u0441\u043a\u0438\u0439\" line.\n#\n# The target text (\" \u041c. \u0421. \u0410\u043d\u0430\u043d\u044c\u0435\u0432\u0441\u043a\u0438\u0439\") is, in the canonical edit, split\n# into six sibling runs instead of being collapsed into one:\n#   \" \" | \"\u0415\" | \". \" | \"\u041c\" | \". \" | \"\u041a\u0443\u0437\u043d\u0435\u0446\u043e\u0432\u0430\"\n# A plain Range.Text (or Find.Execute with Replacement.Text) assignment\n# would leave the new text in a single run (and this host even coalesces it\n# into the preceding \"___\" run, since both share identical formatting), so\n# we rebuild the exact run sequence with Range.InsertXML, which inserts the\n# literal OOXML runs verbatim without merging them.\n\n$d = $word.ActiveDocument\n\n# Locate the target text; Find.Execute narrows the range it's called on to\n# the match, so read the match bounds off of it.\n$seek = $d.Content\n$seek.Find.Text = \" \u041c. \u0421. \u0410\u043d\u0430\u043d\u044c\u0435\u0432\u0441\u043a\u0438\u0439\"\n$found = $seek.Find.Execute()\nif (-not $found) {\n  throw \"Target text ' \u041c. \u0421. \u0410\u043d\u0430\u043d\u044c\u0435\u0432\u0441\u043a\u0438\u0439' was not found in the document.\"\n}\n\n# Re-seat a fresh Range over the exact same character span. InsertXML on a\n# Range that's still wired to a just-executed Find object only appends; a\n# plain Document.Range(start, end) replaces the span's contents as intended.\n$target = $d.Range($seek.Start, $seek.End)\n\n$xml = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\"><pkg:xmlData>\n<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\"><Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/></Relationships>\n</pkg:xmlData></pkg:part>\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:rPr><w:rFonts w:cs=\"Times New Roman\"/><w:szCs w:val=\"28\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:rPr><w:rFonts w:cs=\"Times New Roman\"/><w:szCs w:val=\"28\"/></w:rPr><w:t>\u0415</w:t></w:r><w:r><w:rPr><w:rFonts w:cs=\"Times New Roman\"/><w:szCs w:val=\"28\"/></w:rPr><w:t xml:space=\"preserve\">. </w:t></w:r><w:r><w:rPr><w:rFonts w:cs=\"Times New Roman\"/><w:szCs w:val=\"28\"/></w:rPr><w:t>\u041c</w:t></w:r><w:r><w:rPr><w:rFonts w:cs=\"Times New Roman\"/><w:szCs w:val=\"28\"/></w:rPr><w:t xml:space=\"preserve\">. </w:t></w:r><w:r><w:rPr><w:rFonts w:cs=\"Times New Roman\"/><w:szCs w:val=\"28\"/></w:rPr><w:t>\u041a\u0443\u0437\u043d\u0435\u0446\u043e\u0432\u0430</w:t></w:r></w:p></w:body></w:document>\n</pkg:xmlData></pkg:part>\n</pkg:package>\n'@\n\n$target.InsertXML($xml)\n"}
